$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.8997134670487106
$ws.Range("B3").Value = 0.8826291079812206
$ws.Range("B4").Value = 0.9494949494949495
$ws.Range("B5").Value = 0.9148418491484185
